$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.46%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'31.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.65%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.143"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.89%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07369"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.06%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'2.564"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'69.07%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'7.909"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.74%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.756"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.43%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.9183"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.21%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1740"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.66%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07482"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-6.26%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.08134"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.95%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.03045"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.57%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.09925"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.26%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.001511"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.87%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.006165"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.12%"
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'-0.50%"
$ws.Range("E17").Style = "Normal"

$ws.Range("D19").Value = "'0.3290"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.03%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1340"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.74%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'4.657"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'3.66%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.04650"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.73%"
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'-2.30%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001225"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.21%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004472"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.79%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001302"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.14%"
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'6.89%"
$ws.Range("E27").Style = "Normal"

$ws.Range("D39").Value = "'0.01722"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-2.31%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.04537"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.81%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007121"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.63%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1345"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.17%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002214"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.89%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.01093"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-14.61%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00006296"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.43%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.01002"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-23.04%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'1.847"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'160.46%"
$ws.Range("E47").Style = "Normal"
